$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.60"
$ws.Range("D3").Value = "'22.66"
$ws.Range("D4").Value = "'5.389"
$ws.Range("D5").Value = "'0.05696"
$ws.Range("D6").Value = "'3.406"
$ws.Range("D7").Value = "'6.320"
$ws.Range("D8").Value = "'0.8124"
$ws.Range("D9").Value = "'0.9237"
$ws.Range("D10").Value = "'0.0005849"
$ws.Range("D11").Value = "'0.1411"
$ws.Range("D12").Value = "'0.07433"
$ws.Range("D13").Value = "'0.03080"
$ws.Range("D14").Value = "'0.03016"
$ws.Range("D15").Value = "'0.09381"
$ws.Range("D16").Value = "'3.744"
$ws.Range("D17").Value = "'0.001572"
$ws.Range("D18").Value = "'0.04742"
$ws.Range("D19").Value = "'0.01827"
$ws.Range("D20").Value = "'0.006458"
$ws.Range("D21").Value = "'0.005000"
$ws.Range("D22").Value = "'0.001023"
$ws.Range("D23").Value = "'0.0001500"
$ws.Range("D24").Value = "'3.700"
$ws.Range("D26").Value = "'0.3254"
$ws.Range("D40").Value = "'0.03988"
$ws.Range("D41").Value = "'0.006839"
$ws.Range("D43").Value = "'0.002710"
$ws.Range("D44").Value = "'0.007512"
$ws.Range("D45").Value = "'0.00005800"
$ws.Range("D47").Value = "'0.4299"
$ws.Range("D48").Value = "'0.2153"
